# Refresh Universalis market-price snapshots and recomputed Leve profit
# columns (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the crafting-class sheets. Column layout (all sheets):
#   H=currentAveragePrice  I=currentAveragePriceNQ  J=currentAveragePriceHQ
#   K=LevePriceNQ          L=LevePriceHQ            M=LeveProfitNQ  N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    ,@(132, 8, 20081500)
    ,@(132, 9, 21362944)
    ,@(132, 11, 64088832)
    ,@(132, 13, -64086302)
    ,@(137, 8, 3978.2827)
    ,@(137, 9, 4040)
    ,@(137, 10, 3890.5789)
    ,@(137, 11, 12120)
    ,@(137, 12, 11671.7367)
    ,@(137, 13, -9570)
    ,@(137, 14, -16771.7367)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    ,@(61, 8, 2660.2666)
    ,@(61, 9, 1801.4286)
    ,@(61, 10, 3411.75)
    ,@(61, 11, 1801.4286)
    ,@(61, 12, 3411.75)
    ,@(61, 13, -1589.4286)
    ,@(61, 14, -3835.75)
    ,@(63, 8, 10658746)
    ,@(63, 9, 27704242)
    ,@(63, 10, 5310)
    ,@(63, 11, 27704242)
    ,@(63, 12, 5310)
    ,@(63, 13, -27703556)
    ,@(63, 14, -6682)
    ,@(66, 8, 10658746)
    ,@(66, 9, 27704242)
    ,@(66, 10, 5310)
    ,@(66, 11, 138521210)
    ,@(66, 12, 26550)
    ,@(66, 13, -138517778)
    ,@(66, 14, -33414)
    ,@(74, 8, 3985.5676)
    ,@(74, 9, 4452.154)
    ,@(74, 11, 4452.154)
    ,@(74, 13, -3578.154)
    ,@(77, 8, 3985.5676)
    ,@(77, 9, 4452.154)
    ,@(77, 11, 22260.77)
    ,@(77, 13, -17892.77)
    ,@(122, 8, 4806.4287)
    ,@(122, 9, 2004)
    ,@(122, 10, 6908.25)
    ,@(122, 11, 6012)
    ,@(122, 12, 20724.75)
    ,@(122, 13, -3562)
    ,@(122, 14, -25624.75)
    ,@(132, 8, 3121.182)
    ,@(132, 9, 1584.8334)
    ,@(132, 10, 4964.8)
    ,@(132, 11, 4754.5002)
    ,@(132, 12, 14894.4)
    ,@(132, 13, -2224.5002)
    ,@(132, 14, -19954.4)
    ,@(136, 8, 2660.2666)
    ,@(136, 9, 1801.4286)
    ,@(136, 10, 3411.75)
    ,@(136, 11, 5404.2858)
    ,@(136, 12, 10235.25)
    ,@(136, 13, -2854.2858)
    ,@(136, 14, -15335.25)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    ,@(137, 8, 56922.668)
    ,@(137, 10, 56922.668)
    ,@(137, 12, 56922.668)
    ,@(137, 14, -67122.66800000001)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    ,@(31, 8, 7354.5454)
    ,@(31, 9, 0)
    ,@(31, 10, 7354.5454)
    ,@(31, 11, 0)
    ,@(31, 12, 7354.5454)
    ,@(31, 14, -7944.5454)
    ,@(34, 8, 7354.5454)
    ,@(34, 9, 0)
    ,@(34, 10, 7354.5454)
    ,@(34, 11, 0)
    ,@(34, 12, 7354.5454)
    ,@(34, 14, -7758.5454)
    ,@(58, 8, 2311.1167)
    ,@(58, 9, 1810.5186)
    ,@(58, 11, 1810.5186)
    ,@(58, 13, -1607.5186)
    ,@(132, 8, 4973.316)
    ,@(132, 9, 4593.857)
    ,@(132, 11, 13781.571)
    ,@(132, 13, -11251.571)
    ,@(136, 8, 2311.1167)
    ,@(136, 9, 1810.5186)
    ,@(136, 11, 5431.5558)
    ,@(136, 13, -2881.5558)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    ,@(2, 8, 4464376)
    ,@(2, 9, 105)
    ,@(2, 11, 630)
    ,@(2, 13, -517)
    ,@(5, 8, 2225.5557)
    ,@(5, 10, 3954.6155)
    ,@(5, 12, 11863.8465)
    ,@(5, 14, -12087.8465)
    ,@(23, 8, 211.23529)
    ,@(23, 9, 109.25)
    ,@(23, 10, 242.61539)
    ,@(23, 11, 327.75)
    ,@(23, 12, 727.84617)
    ,@(23, 13, -92.75)
    ,@(23, 14, -1197.84617)
    ,@(75, 8, 1081.75)
    ,@(75, 10, 1038.3334)
    ,@(75, 12, 3115.0002)
    ,@(75, 14, -5111.0002)
    ,@(78, 8, 1081.75)
    ,@(78, 10, 1038.3334)
    ,@(78, 12, 9345.000599999999)
    ,@(78, 14, -19329.0006)
    ,@(113, 8, 596.09753)
    ,@(113, 9, 594.36)
    ,@(113, 10, 598.8125)
    ,@(113, 11, 1783.08)
    ,@(113, 12, 1796.4375)
    ,@(113, 13, 386.9200000000001)
    ,@(113, 14, -6136.4375)
    ,@(129, 8, 2939.261)
    ,@(129, 10, 3314.7778)
    ,@(129, 12, 9944.3334)
    ,@(129, 14, -19944.3334)
    ,@(132, 8, 1833.3889)
    ,@(132, 9, 1160.9231)
    ,@(132, 10, 2213.4783)
    ,@(132, 11, 10448.3079)
    ,@(132, 12, 19921.3047)
    ,@(132, 13, -7918.3079)
    ,@(132, 14, -24981.3047)
    ,@(133, 8, 3965.5557)
    ,@(133, 9, 4268.8887)
    ,@(133, 10, 3662.2222)
    ,@(133, 11, 12806.6661)
    ,@(133, 12, 10986.6666)
    ,@(133, 13, -7746.666100000002)
    ,@(133, 14, -21106.6666)
    ,@(135, 8, 2225.5557)
    ,@(135, 10, 3954.6155)
    ,@(135, 12, 35591.5395)
    ,@(135, 14, -40661.5395)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    ,@(7, 8, 5025.2144)
    ,@(7, 9, 2413)
    ,@(7, 10, 8508.166999999999)
    ,@(7, 11, 2413)
    ,@(7, 12, 8508.166999999999)
    ,@(7, 13, -2301)
    ,@(7, 14, -8732.166999999999)
    ,@(87, 8, 25085.5)
    ,@(87, 9, 10171)
    ,@(87, 10, 40000)
    ,@(87, 11, 10171)
    ,@(87, 12, 40000)
    ,@(87, 13, -9048)
    ,@(87, 14, -42246)
    ,@(90, 8, 25085.5)
    ,@(90, 9, 10171)
    ,@(90, 10, 40000)
    ,@(90, 11, 30513)
    ,@(90, 12, 120000)
    ,@(90, 13, -24897)
    ,@(90, 14, -131232)
    ,@(122, 8, 3756.8438)
    ,@(122, 9, 2513.5908)
    ,@(122, 11, 7540.7724)
    ,@(122, 13, -5090.7724)
    ,@(126, 8, 5025.2144)
    ,@(126, 9, 2413)
    ,@(126, 10, 8508.166999999999)
    ,@(126, 11, 7239)
    ,@(126, 12, 25524.501)
    ,@(126, 13, -4769)
    ,@(126, 14, -30464.501)
    ,@(132, 8, 4604.0786)
    ,@(132, 9, 2114.1667)
    ,@(132, 10, 8161.095)
    ,@(132, 11, 6342.500100000001)
    ,@(132, 12, 24483.285)
    ,@(132, 13, -3812.500100000001)
    ,@(132, 14, -29543.285)
    ,@(136, 8, 4012.15)
    ,@(136, 9, 1745.375)
    ,@(136, 10, 5523.3335)
    ,@(136, 11, 5236.125)
    ,@(136, 12, 16570.0005)
    ,@(136, 13, -2686.125)
    ,@(136, 14, -21670.0005)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    ,@(126, 8, 509765.75)
    ,@(126, 9, 1898.6666)
    ,@(126, 10, 1186921.9)
    ,@(126, 11, 5695.9998)
    ,@(126, 12, 3560765.7)
    ,@(126, 13, -3225.9998)
    ,@(126, 14, -3565705.7)
    ,@(136, 8, 6818)
    ,@(136, 9, 1665.5)
    ,@(136, 11, 4996.5)
    ,@(136, 13, -2446.5)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 13).ClearContents()   # M31 removed (merged into recomputed L31)
$ws.Cells.Item(34, 13).ClearContents()   # M34 removed (merged into recomputed L34)

